$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profile")

# Update the "Profile" sheet sample row to reflect the new Edit Profile
# test options (AvailableTime / Hours / EarnTarget answers).
$ws.Range("C2").Value = "Between $500 and $1000 per month"
$ws.Range("D2").Value = "Full Time"
$ws.Range("E2").Value = "As needed"

# Move the active selection to E10, matching the updated test state.
$ws.Range("E10").Select()
